# Apply "feat: add 2022-Q3 data" edit:
#  - total/"总计" sheet: add a new row for the 2022-Q3 summary (2 funds, 0.05 billion)
#    and push the existing 2022-Q2 / 2022-Q1 summary rows down.
#  - insert a brand-new "2022-Q3" worksheet (fund-level detail) right after "总计",
#    ahead of the existing "2022-Q2"/"2022-Q1" worksheets.
#  - the old "2022-Q2" and "2022-Q1" worksheets keep their data untouched, they are
#    just recreated (same content) so that the sheetId/tab order line up.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Helper: write a "numeric-looking" value as genuine TEXT (not a Number) the
# way the original workbook stores fund codes / ratios, e.g. "000520", "5.81".
# Plain `.Value = "5.81"` gets auto-coerced to a number by the engine, and
# `.Value = "000520"` loses its leading zeroes - so cells holding such values
# are pre-formatted as Text ("@") before the value is poked in.
# ---------------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ===========================================================================
# 1) "总计" sheet (sheetId=1, untouched identity) - rewrite the summary rows.
# ===========================================================================
$total = $wb.Worksheets.Item(1)

# Row 2 becomes the new 2022-Q3 summary (was 2022-Q2).
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.05

# Row 3 becomes the 2022-Q2 summary (was 2022-Q1); A3 already equals 1.
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.21

# Row 4 is brand new: the 2022-Q1 summary. Stamp column A with the same style
# used by A2/A3 (copy format only) before filling in the values.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial($xlPasteFormats)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.35

$total.Range("A1").Select()

# ===========================================================================
# 2) Recreate the per-quarter detail sheets so sheetId or tab order end up
#    exactly as: 总计(1), 2022-Q3(2), 2022-Q2(3), 2022-Q1(4).
# ===========================================================================

# Pin a swatch of the "2022-Q1" header style (cellXfs index used only there)
# on the "总计" sheet so it survives deleting the "2022-Q1" worksheet.
$q1Old = $wb.Worksheets.Item("2022-Q1")
$q1Old.Range("B1").Copy()
$total.Range("ZZ1").PasteSpecial($xlPasteFormats)

# A swatch of the "2022-Q2" / "总计" header style (cellXfs index 2).
$total.Range("B1").Copy()
$total.Range("ZZ2").PasteSpecial($xlPasteFormats)

$wb.Worksheets.Item("2022-Q2").Delete()
$wb.Worksheets.Item("2022-Q1").Delete()

# ---------------------------------------------------------------------------
# 2022-Q3 (brand new data, style like the old "2022-Q2" sheet)
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$total.Range("ZZ2").Copy()
$q3.Range("B1:H1").PasteSpecial($xlPasteFormats)
$q3.Range("A2:A3").PasteSpecial($xlPasteFormats)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "004138"
$q3.Range("C2").Value = "上银鑫达灵活配置混合A"
Set-TextValue $q3.Range("D2") "1.86"
Set-TextValue $q3.Range("E2") "75.90"
Set-TextValue $q3.Range("F2") "2.83"
Set-TextValue $q3.Range("G2") "0.0526"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "015753"
$q3.Range("C3").Value = "上银鑫达灵活配置混合C"
Set-TextValue $q3.Range("D3") "0.03"
Set-TextValue $q3.Range("E3") "75.90"
Set-TextValue $q3.Range("F3") "2.83"
Set-TextValue $q3.Range("G3") "0.0008"
$q3.Range("H3").Value = 10

# Reset the format of the text-like cells back to the plain (unstyled) look
# used in the source file (NumberFormat="@" above only exists to stop the
# engine from mangling the values - the actual cells carry no explicit style).
$total.Range("C2").Copy()
$q3.Range("B2:G3").PasteSpecial($xlPasteFormats)
$q3.Range("A2:A3").PasteSpecial($xlPasteFormats) # restore col-A style after the blanket reset above
$total.Range("ZZ2").Copy()
$q3.Range("A2:A3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 2022-Q2 (unchanged data, same style as before)
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Add($null, $q3)
$q2.Name = "2022-Q2"

$total.Range("ZZ2").Copy()
$q2.Range("B1:H1").PasteSpecial($xlPasteFormats)
$q2.Range("A2:A3").PasteSpecial($xlPasteFormats)

$q2.Range("B1").Value = "基金代码"
$q2.Range("C1").Value = "基金名称"
$q2.Range("D1").Value = "基金规模"
$q2.Range("E1").Value = "股票总仓位"
$q2.Range("F1").Value = "仓位占比"
$q2.Range("G1").Value = "持有市值(亿元)"
$q2.Range("H1").Value = "仓位排名"

$q2.Range("A2").Value = 0
Set-TextValue $q2.Range("B2") "000520"
$q2.Range("C2").Value = "上银新兴价值成长混合"
Set-TextValue $q2.Range("D2") "5.81"
Set-TextValue $q2.Range("E2") "77.13"
Set-TextValue $q2.Range("F2") "2.48"
Set-TextValue $q2.Range("G2") "0.1441"
$q2.Range("H2").Value = 9

$q2.Range("A3").Value = 1
Set-TextValue $q2.Range("B3") "004138"
$q2.Range("C3").Value = "上银鑫达灵活配置混合"
Set-TextValue $q2.Range("D3") "2.22"
Set-TextValue $q2.Range("E3") "79.09"
Set-TextValue $q2.Range("F3") "2.88"
Set-TextValue $q2.Range("G3") "0.0639"
$q2.Range("H3").Value = 9

$total.Range("C2").Copy()
$q2.Range("B2:G3").PasteSpecial($xlPasteFormats)
$total.Range("ZZ2").Copy()
$q2.Range("A2:A3").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 2022-Q1 (unchanged data, style pinned from the deleted original sheet)
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q2)
$q1.Name = "2022-Q1"

$total.Range("ZZ1").Copy()
$q1.Range("B1:H1").PasteSpecial($xlPasteFormats)
$q1.Range("A2:A5").PasteSpecial($xlPasteFormats)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "004350"
$q1.Range("C2").Value = "汇丰晋信价值先锋股票"
Set-TextValue $q1.Range("D2") "4.99"
Set-TextValue $q1.Range("E2") "93.32"
Set-TextValue $q1.Range("F2") "2.75"
Set-TextValue $q1.Range("G2") "0.1372"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "000520"
$q1.Range("C3").Value = "上银新兴价值成长混合"
Set-TextValue $q1.Range("D3") "4.64"
Set-TextValue $q1.Range("E3") "77.73"
Set-TextValue $q1.Range("F3") "2.20"
Set-TextValue $q1.Range("G3") "0.1021"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
Set-TextValue $q1.Range("B4") "004138"
$q1.Range("C4").Value = "上银鑫达灵活配置混合"
Set-TextValue $q1.Range("D4") "2.41"
Set-TextValue $q1.Range("E4") "77.12"
Set-TextValue $q1.Range("F4") "3.30"
Set-TextValue $q1.Range("G4") "0.0795"
$q1.Range("H4").Value = 9

$q1.Range("A5").Value = 3
Set-TextValue $q1.Range("B5") "007393"
$q1.Range("C5").Value = "上银未来生活灵活配置混合"
Set-TextValue $q1.Range("D5") "1.05"
Set-TextValue $q1.Range("E5") "85.99"
Set-TextValue $q1.Range("F5") "2.96"
Set-TextValue $q1.Range("G5") "0.0311"
$q1.Range("H5").Value = 10

$total.Range("C2").Copy()
$q1.Range("B2:G5").PasteSpecial($xlPasteFormats)
$total.Range("ZZ1").Copy()
$q1.Range("A2:A5").PasteSpecial($xlPasteFormats)

# Page margins on the 2022-Q1 sheet differ from the workbook default.
$q1.PageSetup.LeftMargin = 0.7 * 72
$q1.PageSetup.RightMargin = 0.7 * 72
$q1.PageSetup.TopMargin = 0.75 * 72
$q1.PageSetup.BottomMargin = 0.75 * 72
$q1.PageSetup.HeaderMargin = 0.3 * 72
$q1.PageSetup.FooterMargin = 0.3 * 72

# Clear out the temporary style-pinning swatch cells on "总计".
$total.Range("ZZ1").Clear()
$total.Range("ZZ2").Clear()

# The original file had the 2022-Q1 tab as the selected one.
$q1.Select()
$total.Range("A1").Select()
$q1.Select()
